$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify the plot_trend_status() calls stored in the "plot" column (F):
# drop the now-redundant parameter="..." argument since it duplicates df=.
$ws.Range("F68").Value = 'plot_trend_status(df=fish_weight, mpa=MPAs, area="Western/Emerald Banks Conservation Area (Restricted Fisheries Zone)", type=NULL)'
$ws.Range("F69").Value = 'plot_trend_status(df=fish_length, mpa=MPAs, area="Western/Emerald Banks Conservation Area (Restricted Fisheries Zone)", type=NULL)'
$ws.Range("F71").Value = 'plot_trend_status(df=zooplankton, mpa=MPAs, type=NULL)'
$ws.Range("F73").Value = 'plot_trend_status(df=haddock_biomass, mpa=MPAs, area="Western/Emerald Banks Conservation Area (Restricted Fisheries Zone)", type=NULL)'
$ws.Range("F74").Value = 'plot_trend_status(df=all_haddock, mpa=MPAs, area="Western/Emerald Banks Conservation Area (Restricted Fisheries Zone)", type=NULL)'
$ws.Range("F76").Value = 'plot_trend_status(df=nitrate, mpa=MPAs, type="surface")'
$ws.Range("F78").Value = 'plot_trend_status(df=salinity, mpa=MPAs, type="surface", dataframe=TRUE)'
$ws.Range("F83").Value = 'plot_trend_status(df=chlorophyll, mpa=MPAs, type="surface")'
$ws.Range("F84").Value = 'plot_trend_status(df=bloom_df, mpa=MPAs, type="surface")'
$ws.Range("F97").Value = 'plot_trend_status(df=temperature, mpa=MPAs, type="surface", dataframe=TRUE)'
$ws.Range("F113").Value = 'plot_trend_status(df=whale_biodiversity, mpa=MPAs, type=NULL)'
$ws.Range("F114").Value = 'plot_trend_status(df=surface_height, mpa=MPAs, type=NULL)'

# Update the saved view state (scroll position / selection) to match the
# author's final position in the sheet when the workbook was saved.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Application.ActiveWindow.ScrollRow = 106
$ws.Range("D118").Select()
